# Update Marano.xlsx sheet data through row 328 (data through 25 July 2021 /
# "aggiornamento fino a 28 luglio" per commit message), appending rows
# 302-328 to the existing A1:D301 table of daily new-positive COVID counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: row, date-serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(302, 44376, 0, 5, 94.6969696969697),
    @(303, 44377, 2, 6, 113.6363636363636),
    @(304, 44378, 0, 5, 94.6969696969697),
    @(305, 44379, 1, 6, 113.6363636363636),
    @(306, 44380, 0, 5, 94.6969696969697),
    @(307, 44381, 0, 4, 75.75757575757575),
    @(308, 44382, 2, 5, 94.6969696969697),
    @(309, 44383, 0, 5, 94.6969696969697),
    @(310, 44384, 0, 3, 56.81818181818181),
    @(311, 44385, 0, 3, 56.81818181818181),
    @(312, 44386, 2, 4, 75.75757575757575),
    @(313, 44387, 0, 4, 75.75757575757575),
    @(314, 44388, 0, 4, 75.75757575757575),
    @(315, 44389, 1, 3, 56.81818181818181),
    @(316, 44390, 0, 3, 56.81818181818181),
    @(317, 44391, 0, 3, 56.81818181818181),
    @(318, 44392, 0, 3, 56.81818181818181),
    @(319, 44393, 0, 1, 18.93939393939394),
    @(320, 44394, 0, 1, 18.93939393939394),
    @(321, 44395, 0, 1, 18.93939393939394),
    @(322, 44396, 0, 0, 0),
    @(323, 44397, 0, 0, 0),
    @(324, 44398, 0, 0, 0),
    @(325, 44399, 0, 0, 0),
    @(326, 44400, 0, 0, 0),
    @(327, 44401, 0, 0, 0),
    @(328, 44402, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]

    # Column A carries the date-style formatting (style index reused from
    # the last existing data row, s="2" / YYYY-MM-DD HH:MM:SS). Copying the
    # format from the row above keeps the same style entry instead of
    # minting a new one.
    $ws.Range("A301").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
